$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder match rows within same-date groups (E:AD), keep A (index) and C/D (league/date) fixed.

$ws.Range("B236").Value = 6941439
$row236 = New-Object "object[,]" 1,26
$row236[0,0] = "Al Wehda Mecca"
$row236[0,1] = "Al Ahli Jeddah"
$row236[0,2] = 1
$row236[0,3] = 1
$row236[0,4] = 0
$row236[0,5] = 0
$row236[0,6] = "D"
$row236[0,7] = 5
$row236[0,8] = 4.333
$row236[0,9] = 1.5
$row236[0,10] = 5.25
$row236[0,11] = 4.5
$row236[0,12] = 1.45
$row236[0,13] = 1.25
$row236[0,14] = 1.775
$row236[0,15] = 2.025
$row236[0,16] = 3
$row236[0,17] = 1.925
$row236[0,18] = 1.875
$row236[0,19] = -1
$row236[0,20] = 3.5
$row236[0,21] = -1
$row236[0,22] = 0.7749999999999999
$row236[0,23] = -1
$row236[0,24] = -1
$row236[0,25] = 0.875
$ws.Range("E236:AD236").Value = $row236

$ws.Range("B237").Value = 6940788
$row237 = New-Object "object[,]" 1,26
$row237[0,0] = "Al Ittihad Jeddah"
$row237[0,1] = "Al Taawon Buraidah"
$row237[0,2] = 0
$row237[0,3] = 0
$row237[0,4] = 0
$row237[0,5] = 0
$row237[0,6] = "D"
$row237[0,7] = 1.727
$row237[0,8] = 3.6
$row237[0,9] = 4.333
$row237[0,10] = 1.6
$row237[0,11] = 3.8
$row237[0,12] = 5
$row237[0,13] = -1
$row237[0,14] = 2
$row237[0,15] = 1.8
$row237[0,16] = 3
$row237[0,17] = 1.95
$row237[0,18] = 1.85
$row237[0,19] = -1
$row237[0,20] = 2.8
$row237[0,21] = -1
$row237[0,22] = -1
$row237[0,23] = 0.8
$row237[0,24] = -1
$row237[0,25] = 0.8500000000000001
$ws.Range("E237:AD237").Value = $row237

$ws.Range("B238").Value = 6941441
$row238 = New-Object "object[,]" 1,26
$row238[0,0] = "Damac FC"
$row238[0,1] = "AlNassr Riyadh"
$row238[0,2] = 0
$row238[0,3] = 1
$row238[0,4] = 0
$row238[0,5] = 0
$row238[0,6] = "A"
$row238[0,7] = 4.333
$row238[0,8] = 4.75
$row238[0,9] = 1.55
$row238[0,10] = 2.9
$row238[0,11] = 3.75
$row238[0,12] = 2.05
$row238[0,13] = 0.25
$row238[0,14] = 1.975
$row238[0,15] = 1.825
$row238[0,16] = 2.75
$row238[0,17] = 1.85
$row238[0,18] = 1.95
$row238[0,19] = -1
$row238[0,20] = -1
$row238[0,21] = 1.05
$row238[0,22] = -1
$row238[0,23] = 0.825
$row238[0,24] = -1
$row238[0,25] = 0.95
$ws.Range("E238:AD238").Value = $row238

$ws.Range("B239").Value = 6941442
$row239 = New-Object "object[,]" 1,26
$row239[0,0] = "Al Khaleej Saihat"
$row239[0,1] = "Al Hilal Riyadh"
$row239[0,2] = 1
$row239[0,3] = 4
$row239[0,4] = 1
$row239[0,5] = 2
$row239[0,6] = "A"
$row239[0,7] = 9.5
$row239[0,8] = 6
$row239[0,9] = 1.222
$row239[0,10] = 8.5
$row239[0,11] = 5.75
$row239[0,12] = 1.25
$row239[0,13] = 1.75
$row239[0,14] = 1.9
$row239[0,15] = 1.9
$row239[0,16] = 3.25
$row239[0,17] = 1.85
$row239[0,18] = 1.95
$row239[0,19] = -1
$row239[0,20] = -1
$row239[0,21] = 0.25
$row239[0,22] = -1
$row239[0,23] = 0.8999999999999999
$row239[0,24] = 0.8500000000000001
$row239[0,25] = -1
$ws.Range("E239:AD239").Value = $row239

$ws.Range("B282").Value = 7118375
$row282 = New-Object "object[,]" 1,26
$row282[0,0] = "Al Wehda Mecca"
$row282[0,1] = "Al Raed"
$row282[0,2] = 0
$row282[0,3] = 1
$row282[0,4] = 0
$row282[0,5] = 0
$row282[0,6] = "A"
$row282[0,7] = 2.45
$row282[0,8] = 3.4
$row282[0,9] = 2.55
$row282[0,10] = 2.5
$row282[0,11] = 3.5
$row282[0,12] = 2.45
$row282[0,13] = 0
$row282[0,14] = 1.925
$row282[0,15] = 1.875
$row282[0,16] = 2.75
$row282[0,17] = 1.95
$row282[0,18] = 1.85
$row282[0,19] = -1
$row282[0,20] = -1
$row282[0,21] = 1.45
$row282[0,22] = -1
$row282[0,23] = 0.875
$row282[0,24] = -1
$row282[0,25] = 0.8500000000000001
$ws.Range("E282:AD282").Value = $row282

$ws.Range("B283").Value = 7118421
$row283 = New-Object "object[,]" 1,26
$row283[0,0] = "Al Khaleej Saihat"
$row283[0,1] = "Al Ittihad Jeddah"
$row283[0,2] = 1
$row283[0,3] = 1
$row283[0,4] = 0
$row283[0,5] = 1
$row283[0,6] = "D"
$row283[0,7] = 3.3
$row283[0,8] = 3.75
$row283[0,9] = 1.9
$row283[0,10] = 2.75
$row283[0,11] = 3.7
$row283[0,12] = 2.15
$row283[0,13] = 0.25
$row283[0,14] = 1.85
$row283[0,15] = 1.95
$row283[0,16] = 3
$row283[0,17] = 1.975
$row283[0,18] = 1.825
$row283[0,19] = -1
$row283[0,20] = 2.7
$row283[0,21] = -1
$row283[0,22] = 0.425
$row283[0,23] = -0.5
$row283[0,24] = -1
$row283[0,25] = 0.825
$ws.Range("E283:AD283").Value = $row283

$ws.Range("B290").Value = 7138411
$row290 = New-Object "object[,]" 1,26
$row290[0,0] = "Al Raed"
$row290[0,1] = "Al Ahli Jeddah"
$row290[0,2] = 0
$row290[0,3] = 0
$row290[0,4] = 0
$row290[0,5] = 0
$row290[0,6] = "D"
$row290[0,7] = 4.5
$row290[0,8] = 4.2
$row290[0,9] = 1.615
$row290[0,10] = 4
$row290[0,11] = 3.8
$row290[0,12] = 1.75
$row290[0,13] = 0.75
$row290[0,14] = 1.8
$row290[0,15] = 2
$row290[0,16] = 2.75
$row290[0,17] = 1.975
$row290[0,18] = 1.825
$row290[0,19] = -1
$row290[0,20] = 2.8
$row290[0,21] = -1
$row290[0,22] = 0.8
$row290[0,23] = -1
$row290[0,24] = -1
$row290[0,25] = 0.825
$ws.Range("E290:AD290").Value = $row290

$ws.Range("B291").Value = 7138412
$row291 = New-Object "object[,]" 1,26
$row291[0,0] = "Al Akhdoud"
$row291[0,1] = "Al Wehda Mecca"
$row291[0,2] = 1
$row291[0,3] = 1
$row291[0,4] = 1
$row291[0,5] = 1
$row291[0,6] = "D"
$row291[0,7] = 1.909
$row291[0,8] = 3.75
$row291[0,9] = 3.4
$row291[0,10] = 1.615
$row291[0,11] = 4.1
$row291[0,12] = 4.75
$row291[0,13] = -0.75
$row291[0,14] = 1.775
$row291[0,15] = 2.025
$row291[0,16] = 2.75
$row291[0,17] = 1.9
$row291[0,18] = 1.9
$row291[0,19] = -1
$row291[0,20] = 3.1
$row291[0,21] = -1
$row291[0,22] = -1
$row291[0,23] = 1.025
$row291[0,24] = -1
$row291[0,25] = 0.8999999999999999
$ws.Range("E291:AD291").Value = $row291

$ws.Range("B292").Value = 7141178
$row292 = New-Object "object[,]" 1,26
$row292[0,0] = "Al Ittihad Jeddah"
$row292[0,1] = "Damac FC"
$row292[0,2] = 4
$row292[0,3] = 1
$row292[0,4] = 1
$row292[0,5] = 0
$row292[0,6] = "H"
$row292[0,7] = 1.65
$row292[0,8] = 4.333
$row292[0,9] = 4.333
$row292[0,10] = 1.65
$row292[0,11] = 4.5
$row292[0,12] = 4.2
$row292[0,13] = -0.75
$row292[0,14] = 1.775
$row292[0,15] = 2.025
$row292[0,16] = 3.25
$row292[0,17] = 1.9
$row292[0,18] = 1.9
$row292[0,19] = 0.6499999999999999
$row292[0,20] = -1
$row292[0,21] = -1
$row292[0,22] = 0.7749999999999999
$row292[0,23] = -1
$row292[0,24] = 0.8999999999999999
$row292[0,25] = -1
$ws.Range("E292:AD292").Value = $row292

$ws.Range("B293").Value = 7138410
$row293 = New-Object "object[,]" 1,26
$row293[0,0] = "Al Ittifaq Dammam"
$row293[0,1] = "Al Shabab Riyadh"
$row293[0,2] = 1
$row293[0,3] = 0
$row293[0,4] = 1
$row293[0,5] = 0
$row293[0,6] = "H"
$row293[0,7] = 2.3
$row293[0,8] = 3.5
$row293[0,9] = 2.8
$row293[0,10] = 3.1
$row293[0,11] = 3.6
$row293[0,12] = 2.15
$row293[0,13] = 0.25
$row293[0,14] = 1.95
$row293[0,15] = 1.85
$row293[0,16] = 2.75
$row293[0,17] = 1.95
$row293[0,18] = 1.85
$row293[0,19] = 2.1
$row293[0,20] = -1
$row293[0,21] = -1
$row293[0,22] = 0.95
$row293[0,23] = -1
$row293[0,24] = -1
$row293[0,25] = 0.8500000000000001
$ws.Range("E293:AD293").Value = $row293

$ws.Range("B294").Value = 7141179
$row294 = New-Object "object[,]" 1,26
$row294[0,0] = "Al Fayha"
$row294[0,1] = "Al Taawon Buraidah"
$row294[0,2] = 1
$row294[0,3] = 1
$row294[0,4] = 1
$row294[0,5] = 1
$row294[0,6] = "D"
$row294[0,7] = 3
$row294[0,8] = 3.5
$row294[0,9] = 2.2
$row294[0,10] = 4.333
$row294[0,11] = 4
$row294[0,12] = 1.7
$row294[0,13] = 0.75
$row294[0,14] = 1.95
$row294[0,15] = 1.85
$row294[0,16] = 3
$row294[0,17] = 1.975
$row294[0,18] = 1.825
$row294[0,19] = -1
$row294[0,20] = 3
$row294[0,21] = -1
$row294[0,22] = 0.95
$row294[0,23] = -1
$row294[0,24] = -1
$row294[0,25] = 0.825
$ws.Range("E294:AD294").Value = $row294

$ws.Range("B296").Value = 7141394
$row296 = New-Object "object[,]" 1,26
$row296[0,0] = "Al Hilal Riyadh"
$row296[0,1] = "Al Taee"
$row296[0,2] = 3
$row296[0,3] = 1
$row296[0,4] = 2
$row296[0,5] = 0
$row296[0,6] = "H"
$row296[0,7] = 1.142
$row296[0,8] = 8.5
$row296[0,9] = 11
$row296[0,10] = 1.111
$row296[0,11] = 11
$row296[0,12] = 13
$row296[0,13] = -2.75
$row296[0,14] = 1.875
$row296[0,15] = 1.925
$row296[0,16] = 4.25
$row296[0,17] = 1.825
$row296[0,18] = 1.975
$row296[0,19] = 0.111
$row296[0,20] = -1
$row296[0,21] = -1
$row296[0,22] = -1
$row296[0,23] = 0.925
$row296[0,24] = -0.5
$row296[0,25] = 0.4875
$ws.Range("E296:AD296").Value = $row296

$ws.Range("B297").Value = 7141180
$row297 = New-Object "object[,]" 1,26
$row297[0,0] = "Al Fateh SC"
$row297[0,1] = "Al Hazm"
$row297[0,2] = 2
$row297[0,3] = 1
$row297[0,4] = 2
$row297[0,5] = 1
$row297[0,6] = "H"
$row297[0,7] = 1.7
$row297[0,8] = 4.333
$row297[0,9] = 3.8
$row297[0,10] = 1.444
$row297[0,11] = 5
$row297[0,12] = 5
$row297[0,13] = -1.25
$row297[0,14] = 1.95
$row297[0,15] = 1.85
$row297[0,16] = 3.25
$row297[0,17] = 1.8
$row297[0,18] = 2
$row297[0,19] = 0.444
$row297[0,20] = -1
$row297[0,21] = -1
$row297[0,22] = -0.5
$row297[0,23] = 0.425
$row297[0,24] = -0.5
$row297[0,25] = 0.5
$ws.Range("E297:AD297").Value = $row297

$ws.Range("B299").Value = 7158944
$row299 = New-Object "object[,]" 1,26
$row299[0,0] = "Al Taawon Buraidah"
$row299[0,1] = "Al Ittifaq Dammam"
$row299[0,2] = 1
$row299[0,3] = 0
$row299[0,4] = $null
$row299[0,5] = $null
$row299[0,6] = "H"
$row299[0,7] = 1.909
$row299[0,8] = 3.5
$row299[0,9] = 3.7
$row299[0,10] = 2.05
$row299[0,11] = 3.4
$row299[0,12] = 3.3
$row299[0,13] = -0.25
$row299[0,14] = 1.8
$row299[0,15] = 2
$row299[0,16] = 2.75
$row299[0,17] = 1.95
$row299[0,18] = 1.85
$row299[0,19] = 1.05
$row299[0,20] = -1
$row299[0,21] = -1
$row299[0,22] = 0.8
$row299[0,23] = -1
$row299[0,24] = -1
$row299[0,25] = 0.8500000000000001
$ws.Range("E299:AD299").Value = $row299

$ws.Range("B300").Value = 7153884
$row300 = New-Object "object[,]" 1,26
$row300[0,0] = "Damac FC"
$row300[0,1] = "Al Raed"
$row300[0,2] = 1
$row300[0,3] = 1
$row300[0,4] = $null
$row300[0,5] = $null
$row300[0,6] = "D"
$row300[0,7] = 2.3
$row300[0,8] = 3.4
$row300[0,9] = 2.875
$row300[0,10] = 2.2
$row300[0,11] = 3.5
$row300[0,12] = 3
$row300[0,13] = -0.25
$row300[0,14] = 1.925
$row300[0,15] = 1.875
$row300[0,16] = 2.5
$row300[0,17] = 1.825
$row300[0,18] = 1.975
$row300[0,19] = -1
$row300[0,20] = 2.5
$row300[0,21] = -1
$row300[0,22] = -0.5
$row300[0,23] = 0.4375
$row300[0,24] = -1
$row300[0,25] = 0.9750000000000001
$ws.Range("E300:AD300").Value = $row300

$ws.Range("B301").Value = 7154704
$row301 = New-Object "object[,]" 1,26
$row301[0,0] = "Al Khaleej Saihat"
$row301[0,1] = "Al Riyadh"
$row301[0,2] = 1
$row301[0,3] = 2
$row301[0,4] = $null
$row301[0,5] = $null
$row301[0,6] = "A"
$row301[0,7] = 4
$row301[0,8] = 3.8
$row301[0,9] = 1.75
$row301[0,10] = 2.25
$row301[0,11] = 3.5
$row301[0,12] = 3
$row301[0,13] = -0.25
$row301[0,14] = 1.975
$row301[0,15] = 1.825
$row301[0,16] = 2.75
$row301[0,17] = 2
$row301[0,18] = 1.8
$row301[0,19] = -1
$row301[0,20] = -1
$row301[0,21] = 2
$row301[0,22] = -1
$row301[0,23] = 0.825
$row301[0,24] = 0.5
$row301[0,25] = -0.5
$ws.Range("E301:AD301").Value = $row301

$ws.Range("B302").Value = 7153883
$row302 = New-Object "object[,]" 1,26
$row302[0,0] = "Al Hazm"
$row302[0,1] = "Abha"
$row302[0,2] = 2
$row302[0,3] = 1
$row302[0,4] = $null
$row302[0,5] = $null
$row302[0,6] = "H"
$row302[0,7] = 4.2
$row302[0,8] = 4.2
$row302[0,9] = 1.65
$row302[0,10] = 2.5
$row302[0,11] = 3.9
$row302[0,12] = 2.35
$row302[0,13] = 0
$row302[0,14] = 1.925
$row302[0,15] = 1.875
$row302[0,16] = 3.25
$row302[0,17] = 1.9
$row302[0,18] = 1.9
$row302[0,19] = 1.5
$row302[0,20] = -1
$row302[0,21] = -1
$row302[0,22] = 0.925
$row302[0,23] = -1
$row302[0,24] = -0.5
$row302[0,25] = 0.45
$ws.Range("E302:AD302").Value = $row302

$ws.Range("B304").Value = 7158945
$row304 = New-Object "object[,]" 1,26
$row304[0,0] = "Al Shabab Riyadh"
$row304[0,1] = "Al Fateh SC"
$row304[0,2] = 3
$row304[0,3] = 2
$row304[0,4] = $null
$row304[0,5] = $null
$row304[0,6] = "H"
$row304[0,7] = 1.7
$row304[0,8] = 3.75
$row304[0,9] = 4.5
$row304[0,10] = 1.615
$row304[0,11] = 3.9
$row304[0,12] = 4.75
$row304[0,13] = -1
$row304[0,14] = 1.975
$row304[0,15] = 1.825
$row304[0,16] = 3.25
$row304[0,17] = 1.95
$row304[0,18] = 1.85
$row304[0,19] = 0.615
$row304[0,20] = -1
$row304[0,21] = -1
$row304[0,22] = 0
$row304[0,23] = 0
$row304[0,24] = 0.95
$row304[0,25] = -1
$ws.Range("E304:AD304").Value = $row304
